$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Add quarterly subtotal formulas in column D ---
$ws.Range("D5").Formula  = "=SUM(C2:C5)"
$ws.Range("D9").Formula  = "=SUM(C6:C9)"
$ws.Range("D13").Formula = "=SUM(C10:C13)"
$ws.Range("D17").Formula = "=SUM(C14:C17)"
$ws.Range("D22").Formula = "=SUM(C18:C22)"

# --- Add new row 24 (copy formatting from row 23 first so styles match) ---
$ws.Range("A23:C23").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A24").Value = 41422
$ws.Range("B24").Value = "Besprechung mit Betreuer, Planung anpassen, Versuche mit E3 E4 Mix"
$ws.Range("C24").Value = 3

# --- Update selection to the newly active cell ---
$ws.Range("B24").Select()

# --- Try to update the saved window size/position (best effort) ---
$win = $wb.Windows.Item(1)
$win.Left = -15
$win.Top = -15
$win.Width = 12600
$win.Height = 12405
